$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark (Word's "last edit location" marker)
#    from its old spot to the end of the "Education. Worked with
#    Universities ..." paragraph, as if that paragraph had just been
#    edited.
# ---------------------------------------------------------------------

# Remove the existing (hidden) "_GoBack" bookmark, wherever it is.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the end of the target paragraph's text.
$eduRange = $d.Content
$eduRange.Find.Execute("Education. Worked with Universities around the country. Integrated with external Student systems and Learning Management Systems. Cloud architecture, small and fast paced team. Releases every 15 days, feature-oriented development.")
$eduRange.Collapse(0)

# A genuinely zero-length Range can't be handed straight to Bookmarks.Add
# at a paragraph-end position, so nudge a temporary marker character in,
# bookmark that, then delete the marker -- leaving a clean collapsed
# bookmark exactly where the edit happened.
$eduRange.InsertAfter("Z")
$d.Bookmarks.Add("_GoBack", $eduRange)
$eduRange.Text = ""

# ---------------------------------------------------------------------
# 2) Drop the stale <w:lastRenderedPageBreak/> cached in front of the
#    "Software Engineer at JDA Software" run -- touch the run's text so
#    the layout cache for it gets regenerated.
# ---------------------------------------------------------------------
$jdaRange = $d.Content
$jdaRange.Find.Execute("Software Engineer at JDA Software", $true, $false, $false,
                        $false, $false, $true, 1, $false,
                        "Software Engineer at JDA Software", 2)

# ---------------------------------------------------------------------
# 3) Remove one of the two consecutive blank paragraphs that sit right
#    after the last "...SQLServer antlr.org ISO2001" table.
# ---------------------------------------------------------------------
$lastSkillsTable = $d.Tables(10)
$tblEnd = $lastSkillsTable.Range.End
$blankMark = $d.Range($tblEnd, $tblEnd + 1)
$blankMark.Delete()

# ---------------------------------------------------------------------
# 4) Widen the bottom page margin from 270 twips (13.5pt) to 540 twips
#    (27pt).
# ---------------------------------------------------------------------
$d.Sections(1).PageSetup.BottomMargin = 27
